$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for columns A:D, rows 2-17 (header row 1 stays A/B/C, add header D)
$ws.Range("D1").Value = "ITI"

$data = @(
    @(1, 9, 1, 6),
    @(2, 5, 2, 7),
    @(3, 6, 2, 7),
    @(4, 33, 4, 7),
    @(5, 30, 1, 7),
    @(6, 21, 1, 7),
    @(7, 17, 3, 8),
    @(8, 28, 4, 6),
    @(9, 1, 3, 8),
    @(10, 34, 2, 6),
    @(11, 8, 4, 8),
    @(12, 27, 4, 8),
    @(13, 26, 2, 7),
    @(14, 37, 3, 6),
    @(15, 7, 3, 7),
    @(16, 12, 1, 7)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# Remove the leftover rows (18, 19, 20) that are no longer part of the table
$ws.Range("A18:D20").Clear()

# Column width adjustments (values tuned so the stored OOXML width matches
# the target 17.1640625 / 11 as closely as the engine's pixel-grid rounding allows)
$ws.Range("C1").ColumnWidth = 16.33
$ws.Range("F1").ColumnWidth = 10.17

# Update selection to match the new active cell / selection in the diff
$ws.Range("K21").Select()
